# "all styles applies recursively"
#
# 1) The three existing " both" runs (slides 1, 4, 7) get bold added while
#    keeping their existing underline + Calibri formatting.
# 2) A new slide 8 is appended (duplicate of slide 7's "Title and Content"
#    layout) whose content run reads "That`s all, folks!" in bold, underlined,
#    green (00FF00) Calibri 14pt text.

$p = $ppt.ActivePresentation

# --- 1. Add bold to every " both" run (slides 1, 4 and 7) -----------------
foreach ($idx in 1, 4, 7) {
    $slide = $p.Slides.Item($idx)
    $shape = $slide.Shapes.Item(2)
    $tr = $shape.TextFrame.TextRange
    $runCount = $tr.Runs().Count
    for ($i = 1; $i -le $runCount; $i++) {
        $run = $tr.Runs($i, 1)
        if ($run.Text -eq " both") {
            $run.Font.Bold = $true
        }
    }
}

# --- 2. Append new slide 8 ("That`s all, folks!") --------------------------
$lastSlide = $p.Slides.Item($p.Slides.Count)
$newSlides = $lastSlide.Duplicate()
$newSlide = $newSlides.Item(1)

$content = $newSlide.Shapes.Item(2)
$contentRange = $content.TextFrame.TextRange
$contentRange.Text = "That``s all, folks!"
$contentRange.Font.Bold = $true
$contentRange.Font.Color.RGB = 65280
